# Refresh scraped bilibili-huodong listing counts ("想去人数" / "最低票价")
# across the 展览 / 演出 / 全部类型 sheets, and replace three rows in
# 全部类型 (13-15) with newly scraped events, matching the upstream
# gh-pages data regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 244
$ws.Range("F3").Value = 229
$ws.Range("F4").Value = 255
$ws.Range("F5").Value = 2845
$ws.Range("G7").Value = "已售罄"
$ws.Range("F8").Value = 2188
$ws.Range("F9").Value = 298
$ws.Range("F10").Value = 31
$ws.Range("F11").Value = 427
$ws.Range("F13").Value = 2527
$ws.Range("F15").Value = 1302
$ws.Range("F16").Value = 4614
$ws.Range("F18").Value = 4926
$ws.Range("F19").Value = 1488
$ws.Range("F20").Value = 2828
$ws.Range("F21").Value = 3223
$ws.Range("F23").Value = 1517
$ws.Range("F24").Value = 242
$ws.Range("F25").Value = 830
$ws.Range("F27").Value = 273
$ws.Range("F28").Value = 921
$ws.Range("F29").Value = 1695
$ws.Range("F31").Value = 263
$ws.Range("F32").Value = 650
$ws.Range("F33").Value = 152
$ws.Range("F34").Value = 314
$ws.Range("F35").Value = 383

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 24
$ws.Range("F14").Value = 41
$ws.Range("F17").Value = 47

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 244
$ws.Range("F8").Value = 229
$ws.Range("F10").Value = 255
$ws.Range("F11").Value = 2845
$ws.Range("C13").Value = "北京·原神同人嘉年华"
$ws.Range("D13").Value = "高碑店东路超级蜂巢 5G直播基地"
$ws.Range("E13").Value = "2024.05.18 10:00-05.19 17:00"
$ws.Range("F13").Value = 2188
$ws.Range("G13").Value = 70
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83649"
$ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202404/e6A56bW11712025385291.jpeg"
$ws.Range("C14").Value = "北京·漫无止境2.0"
$ws.Range("D14").Value = "北京市朝阳区安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园L1层thebONE潮街"
$ws.Range("F14").Value = 298
$ws.Range("G14").Value = 68
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=84086"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202404/zoDpqSxB1712821141709.jpeg"
# Leading apostrophe forces text (matches source inlineStr "2024-05-19"),
# otherwise Excel auto-coerces the date-shaped string to a date serial.
$ws.Range("B15").Value = "'2024-05-19"
$ws.Range("C15").Value = "北京·BanG Dream! Only LIVE"
$ws.Range("D15").Value = "工体北路瑞士公寓地下一层 良田Loamy Space"
$ws.Range("E15").Value = "2024.05.19 14:00-05.19 18:00"
$ws.Range("F15").Value = 85
$ws.Range("G15").Value = "已售罄"
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=84764"
$ws.Range("I15").Value = "//i2.hdslb.com/bfs/openplatform/202404/WANjgM311713869761955.png"
$ws.Range("F17").Value = 31
$ws.Range("F18").Value = 427
$ws.Range("F21").Value = 2527
$ws.Range("F22").Value = 1302
$ws.Range("F25").Value = 24
$ws.Range("F26").Value = 4614
$ws.Range("F28").Value = 4926
$ws.Range("F29").Value = 1488
$ws.Range("F30").Value = 2828
$ws.Range("F31").Value = 3223
$ws.Range("F33").Value = 41
$ws.Range("F35").Value = 1517
$ws.Range("F37").Value = 242
$ws.Range("F38").Value = 830
$ws.Range("F40").Value = 273
$ws.Range("F41").Value = 921
$ws.Range("F42").Value = 47
$ws.Range("F43").Value = 1695
$ws.Range("F45").Value = 263
$ws.Range("F46").Value = 650
$ws.Range("F47").Value = 152
$ws.Range("F48").Value = 314
$ws.Range("F49").Value = 383
